# "added new assets from matchsync-ig"
#   - Metadata!B7  (Experimental) : false -> true      (keep as literal text, not boolean)
#   - Metadata!B8  (Date)         : 2022-11-02T14:44:07-05:00 -> 2024-02-19T18:37:26-06:00
#   - Metadata!B14 (Case Sensitive): (empty) -> true    (keep as literal text, not boolean)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Writing the bare word true/false into a cell auto-converts it to a Boolean
# value in this engine (same as real Excel's autodetection). The source
# workbook stores these as plain text, so build the text via a formula and
# then collapse the formula down to its literal value with copy/paste-values
# - this keeps the cell's type as Text and leaves its style untouched.
$ws.Range("B7").Formula = "=""true"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B14").Formula = "=""true"""
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)  # xlPasteValues

# Plain date-like text is not reinterpreted, so a direct value assignment is fine.
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"

$excel.CutCopyMode = 0
